$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings (e.g. "241.52") stay stored as text,
# matching the original inline-string cell type, instead of being
# auto-converted to numbers by Excel when the Value is assigned.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '36.399.42'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '1.934.33'
$ws.Range("E3").Value = '  -2.33%  '
$ws.Range("D5").Value = '241.52'
$ws.Range("E5").Value = '  -1.17%  '
$ws.Range("D6").Value = '0.607'
$ws.Range("E6").Value = '  -3.22%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '56.97'
$ws.Range("E8").Value = '  -3.86%  '
$ws.Range("D9").Value = '0.359'
$ws.Range("E9").Value = '  -4.24%  '
$ws.Range("D10").Value = '0.0837'
$ws.Range("E10").Value = '  +1.59%  '
$ws.Range("D11").Value = '0.103'
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").Value = '2.216.59'
$ws.Range("E12").Value = '  -2.36%  '
$ws.Range("D13").Value = '0.803'
$ws.Range("E13").Value = '  -6.81%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '21.05'
$ws.Range("E14").Value = '  -10.51%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '13.41'
$ws.Range("E15").Value = '  -3.72%  '
$ws.Range("D16").Value = '5.14'
$ws.Range("E16").Value = '  -5.54%  '
$ws.Range("D17").Value = '1.950.75'
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("D18").Value = '36.288.90'
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").Value = '68.93'
$ws.Range("E19").Value = '  -1.41%  '
$ws.Range("D20").Value = '0.0₃0861'
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("D21").Value = '227.52'
$ws.Range("E21").Value = '  -2.80%  '
$ws.Range("D22").Value = '4.96'
$ws.Range("E22").Value = '  -6.54%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = '2.34'
$ws.Range("E24").Value = '  -10.09%  '
$ws.Range("D25").Value = '2.27'
$ws.Range("E25").Value = '  -1.68%  '
$ws.Range("E26").Value = '  -7.50%  '
$ws.Range("D27").Value = '160.55'
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("D28").Value = '0.131'
$ws.Range("E28").Value = '  -0.73%  '
$ws.Range("D29").Value = '19.18'
$ws.Range("E29").Value = '  -3.13%  '
$ws.Range("E30").Value = '  -2.11%  '
$ws.Range("E31").Value = '  -6.07%  '
$ws.Range("D32").Value = '4.55'
$ws.Range("E32").Value = '  -6.96%  '
$ws.Range("D33").Value = '0.0628'
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").Value = '4.16'
$ws.Range("E34").Value = '  -5.77%  '
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").Value = '6.08'
$ws.Range("E36").Value = '  -2.03%  '
$ws.Range("D37").Value = '1.78'
$ws.Range("E37").Value = '  -0.55%  '
$ws.Range("D38").Value = '2.12'
$ws.Range("E38").Value = '  -6.09%  '
$ws.Range("D39").Value = '2.99'
$ws.Range("E39").Value = '  -1.36%  '
$ws.Range("D40").Value = '0.0968'
$ws.Range("E40").Value = '  +0.47%  '
$ws.Range("E41").Value = '  -1.23%  '
$ws.Range("D42").Value = '1.15'
$ws.Range("E42").Value = '  -7.12%  '
$ws.Range("D43").Value = '0.0208'
$ws.Range("E43").Value = '  -2.59%  '
$ws.Range("D44").Value = '15.53'
$ws.Range("E44").Value = '  -4.01%  '
$ws.Range("D45").Value = '1.332.98'
$ws.Range("E45").Value = '  -2.64%  '
$ws.Range("E46").Value = '  -7.09%  '
$ws.Range("D47").Value = '86.30'
$ws.Range("E47").Value = '  -6.46%  '
$ws.Range("D48").Value = '7.06'
$ws.Range("E48").Value = '  -5.66%  '
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("D50").Value = '43.96'
$ws.Range("E50").Value = '  -3.21%  '
$ws.Range("D51").Value = '2.108.40'
$ws.Range("E51").Value = '  -2.45%  '
